# Estatisticas.xlsx - Adicao de resultados KNN para imagens originais
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- Tabela 1 (Sem Pre-processamento) : coluna B = knn -------------------
$ws.Range("B4").Value  = 0.9313
$ws.Range("B5").Value  = 0.9438
$ws.Range("B6").Value  = 0.9688
$ws.Range("B7").Value  = 0.9565
$ws.Range("B8").Value  = 0.9563
$ws.Range("B9").Value  = 0.9625
$ws.Range("B10").Value = 0.9375
$ws.Range("B11").Value = 0.9625
$ws.Range("B12").Value = 0.9813
$ws.Range("B13").Value = 0.9503

# --- Tabela 2 (Com Pre-processamento) : coluna B = knn --------------------
$ws.Range("B18").Value = 0.9255
$ws.Range("B19").Value = 0.8882
$ws.Range("B20").Value = 0.9814
$ws.Range("B21").Value = 0.9748
$ws.Range("B22").Value = 0.9625
$ws.Range("B23").Value = 0.9503
$ws.Range("B24").Value = 0.9241
$ws.Range("B25").Value = 0.9688
$ws.Range("B26").Value = 0.9877
$ws.Range("B27").Value = 0.9563

# --- Remove leftover ROI / roi+bag scratch data in L18:O19 ----------------
$ws.Range("L18:O19").ClearContents()

# --- Update selection to reflect the author's final cursor position -------
$ws.Range("M16").Select()
